$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the spring stiffness inputs (K1-K4) from kN/m to N/m: relabel and
# replace the raw kN/m value with a formula that scales it to N/m (x1000).
$ws.Range("A5").Value = "K1 (N/m)"
$ws.Range("B5").Formula = "=20.5*1000"

$ws.Range("A6").Value = "K2 (N/m)"
$ws.Range("B6").Formula = "=20.5*1000"

$ws.Range("A7").Value = "K3 (N/m)"
$ws.Range("B7").Formula = "=147.2*1000"

$ws.Range("A8").Value = "K4 (N/m)"
$ws.Range("B8").Formula = "=147.2*1000"

# Populate the damping coefficients (C1-C4), previously left at 0, with a
# placeholder excitation/damping value of 100.
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100
$ws.Range("B12").Value = 100

# Torsional stiffness also needs converting from kNm/rad to Nm/rad (x1000).
$ws.Range("B15").Formula = "=12.5*1000"

# Update the active selection left behind on the sheet.
$ws.Range("G15").Select()
